# "switch back to matplotlib" - refresh the rolling 30-day studying window:
# drop the oldest 3 days (2024-03-17/16/15) and add the 3 newest days
# (2024-04-16/15/14) at the top of the "Monthly Data" sheet, and update
# today's completed/remaining hours on the "Today Progress" sheet.

$wb = $excel.ActiveWorkbook

$monthly = $wb.Worksheets.Item("Monthly Data")

# Insert 3 fresh rows right under the header so the newest dates land on top.
$monthly.Rows("2:4").Insert() | Out-Null
# Inserted rows inherit the header's (bold/centered) formatting; strip it so
# the new rows match the plain, unstyled data rows below them.
$monthly.Rows("2:4").ClearFormats() | Out-Null

$newDates  = @("2024-04-16", "2024-04-15", "2024-04-14")
$newHours  = @(4, 1.25, 0)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = 2 + $i

    # Keep the dates as plain text (matching the existing rows) instead of
    # letting Excel auto-convert the "yyyy-mm-dd" string into a date value.
    $dateCell = $monthly.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDates[$i]
    $dateCell.ClearFormats() | Out-Null

    $monthly.Cells.Item($row, 2).Value = $newHours[$i]
}

# The window stays at 30 rows of data, so drop the 3 oldest rows that got
# pushed past the bottom of the range (they're now rows 32:34).
$monthly.Rows("32:34").Delete() | Out-Null

# Update "today's" progress figures.
$today = $wb.Worksheets.Item("Today Progress")
$today.Cells.Item(2, 2).Value = 0.6666666666666666
$today.Cells.Item(3, 2).Value = 3.333333333333333

$excel.Calculate()
